# edit.ps1 - apply the Androidification.docx changes
$d = $word.ActiveDocument

# --- Hunk 1a: split the final run of the "Outlines..." paragraph so that
#     "sure ," is wrapped in a gramStart/gramEnd proofErr pair.
$r = $d.Content
$found = $r.Find.Execute("Outlines no longer work now all", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Outlines no longer work now all' paragraph" }
$p1 = $r.Paragraphs(1)
$target1 = $p1.Range
$xml1 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Outlines no longer work now all </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shaders</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> are allowed, possibly only default was respecting </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>color</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>? Note h p</w:t></w:r><w:r><w:t>h</w:t></w:r><w:r><w:t>ysics is still fine</w:t></w:r><w:r><w:t xml:space="preserve">, I’m using a simple appearance for them so not </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>sure ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> also physics appears to not work either?</w:t></w:r></w:p></pkg:xmlData>
'@
$target1.InsertXML($xml1)

# --- Hunk 1b: insert a brand-new list paragraph right after it.
$r = $d.Content
$found = $r.Find.Execute("also physics appears to not work either?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find end of 'Outlines...' paragraph" }
$pAfter = $r.Paragraphs(1)
$rAfter = $pAfter.Range
$rAfter.Collapse(0)
$rAfter.InsertParagraphAfter()
$rAfter.Collapse(0)
$rAfter.MoveStart(1, 1) | Out-Null
$pNew = $rAfter.Paragraphs(1)
$target2 = $pNew.Range
$xml2 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Newt </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GlWindow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> used and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>allTextureUnitStates.size</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>()==1</w:t></w:r><w:r><w:t xml:space="preserve"> for a given LAND record displays using some weird normal map</w:t></w:r></w:p></pkg:xmlData>
'@
$target2.InsertXML($xml2)

# --- Hunk 2: drop the stray gramStart/gramEnd proofErr wrapping
#     "getBestConfiguration" (spellStart/spellEnd stay).
$r = $d.Content
$found = $r.Find.Execute("getBestConfiguration", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'getBestConfiguration' paragraph" }
$p3 = $r.Paragraphs(1)
$target3 = $p3.Range
$xml3 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>getBestConfiguration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the root of much AWT however possibly new pre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>GLDrawable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> code for Newt dumps this</w:t></w:r></w:p></pkg:xmlData>
'@
$target3.InsertXML($xml3)

# --- Hunk 3: drop the gramStart/gramEnd proofErr wrapping the whole
#     "Lawyer from beth?" sentence.
$r = $d.Content
$found = $r.Find.Execute("Lawyer from", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Lawyer from beth?' paragraph" }
$p4 = $r.Paragraphs(1)
$target4 = $p4.Range
$xml4 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Lawyer from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>beth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>?</w:t></w:r></w:p></pkg:xmlData>
'@
$target4.InsertXML($xml4)
